$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap A3/E3 with A7/E7, and A4/E4 with A6/E6
$ws.Range("A3").Value = 93
$ws.Range("E3").Value = 79.81

$ws.Range("A4").Value = 95
$ws.Range("E4").Value = 79.16

$ws.Range("A6").Value = 97
$ws.Range("E6").Value = 78.89

$ws.Range("A7").Value = 98
$ws.Range("E7").Value = 78.44
